$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data row (2026/01/08, 木, 13, 201) was added to the daily log,
# pushing the existing "2026/12/29 ..." onward block down by one row
# (old row 592 -> new row 593, ..., old row 633 -> new row 634).
$ws.Rows(592).Insert()

# Column A holds a date formatted as plain text (e.g. "2026/12/29"), not a
# real date value. Pre-formatting the cell as Text before assigning keeps
# Excel's automatic date-recognition from converting the literal string
# into a date serial, and re-applying the "Normal" style afterwards drops
# the Text number-format override again so the cell matches its neighbors
# (no explicit style), while keeping the text content intact.
$ws.Range("A592").NumberFormat = "@"
$ws.Range("A592").Value = "2026/01/08"
$ws.Range("A592").Style = "Normal"

$ws.Range("B592").Value = "木"
$ws.Range("C592").Value = 13
$ws.Range("D592").Value = 201
